# "removed random tab from excel"
#
# The workbook originally had 5 tabs: Registration, Login, Parameter,
# random, test_suite. This change deletes the "random" tab, leaving
# Registration, Login, Parameter, test_suite.
#
# It also updates one of the sample email values on the Registration
# sheet (the row that used to share a string with the now-deleted
# "random" sheet) and leaves the workbook's view state (selected cell
# on test_suite, active tab) matching where the author ended up after
# making the edit.

$wb = $excel.ActiveWorkbook

# Update the sample e-mail address on the Registration sheet.
$registration = $wb.Worksheets.Item("Registration")
$registration.Range("A2").Value = "vjgp@test.com"

# Remove the "random" worksheet tab entirely.
$random = $wb.Worksheets.Item("random")
$random.Delete()

# Leave the cell selection on test_suite where the author left it.
$testSuite = $wb.Worksheets.Item("test_suite")
$testSuite.Range("K11").Select()

# The Parameter tab ends up the active tab after the deletion.
$wb.Worksheets.Item("Parameter").Activate()
